$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$ws.Range("B20").Value = "customer3061"
$ws.Range("B21").Value = "CompanyName3029"
$ws.Range("B15").Value = "burhani1095"
$ws.Range("B16").Value = "tno400002318"
$ws.Range("B17").Value = "vat390002438"

$ws.Range("B21").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
